# Apply "create new official titles" edit:
# - Update row 2 (office id 803812 -> 803813, new title/translation/pinyin/source)
# - Append four new rows (3-6) with new official title records

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
$ws.Range("A2").Value = 803813
$ws.Range("B2").Value = "中憲大夫"
# C2 keeps its value ("20") but is stored as text in the sheet, leave as-is
$ws.Range("D2").Value = "Grand Master Exemplar"
$ws.Range("E2").Value = "zhong xian da fu"
$ws.Range("F2").Value = "'68002"

# --- New row 3 ---
$ws.Range("A3").Value = 803814
$ws.Range("B3").Value = "奉直大夫"
$ws.Range("C3").Value = "'20"
$ws.Range("D3").Value = "Grand Master for Forthright Service (Hucker)"
$ws.Range("E3").Value = "feng zhi da fu"
$ws.Range("F3").Value = "'4763"

# --- New row 4 ---
$ws.Range("A4").Value = 803815
$ws.Range("B4").Value = "太醫院吏目"
$ws.Range("C4").Value = "'20"
$ws.Range("D4").Value = "Medical Secretary in the Imperial Academy of Medicine"
$ws.Range("E4").Value = "tai yi yuan li mu"
$ws.Range("F4").Value = "'4763"

# --- New row 5 ---
$ws.Range("A5").Value = 803816
$ws.Range("B5").Value = "遊擊將軍"
$ws.Range("C5").Value = "'20"
$ws.Range("D5").Value = "Mobile Corps Commander (Hucker)"
$ws.Range("E5").Value = "you ji jiang jun"
$ws.Range("F5").Value = "'4763"

# --- New row 6 ---
$ws.Range("A6").Value = 803817
$ws.Range("B6").Value = "鎮同知"
$ws.Range("C6").Value = "'20"
$ws.Range("D6").Value = "Vice Prefect of a Town"
$ws.Range("E6").Value = "zhen tong zhi"
$ws.Range("F6").Value = "'68002"

Write-Output "Applied new official title rows"
